$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.407.99"
$ws.Cells.Item(2, 5).Value = "  +0.01%  "
$ws.Cells.Item(3, 4).Value = "1.849.34"
$ws.Cells.Item(3, 5).Value = "  +0.08%  "
$ws.Cells.Item(4, 5).Value = "  +0.07%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "240.59"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.6283"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "  -0.14%  "
$ws.Cells.Item(7, 5).Value = "  +0.04%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.07656"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = "  +0.69%  "
$ws.Cells.Item(9, 5).Value = "  -0.76%  "
$ws.Cells.Item(10, 5).Value = "  +1.44%  "
$ws.Cells.Item(11, 2).Value = "TRON"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07739"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = "  -0.05%  "
$ws.Cells.Item(12, 2).Value = "Polkadot"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "5.042"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = "  +0.77%  "
$ws.Cells.Item(13, 2).Value = "Polygon"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.6808"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = "  +0.16%  "
$ws.Cells.Item(14, 2).Value = "ShibaInu"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.00001075"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = "  -0.83%  "
$ws.Cells.Item(15, 2).Value = "Litecoin"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "83.34"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).Value = "  -0.46%  "
$ws.Cells.Item(16, 2).Value = "Uniswap"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "6.176"
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).Value = "  -0.08%  "
$ws.Cells.Item(17, 2).Value = "WrappedBTC"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(17, 4).Value = "29.417.57"
$ws.Cells.Item(17, 5).Value = "  -0.01%  "
$ws.Cells.Item(18, 2).Value = "BitcoinCash"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "228.08"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value = "  -0.38%  "
$ws.Cells.Item(19, 2).Value = "Avalanche"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "12.35"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "  -0.82%  "
$ws.Cells.Item(20, 2).Value = "Dai"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "1.001"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = "  +0.06%  "
$ws.Cells.Item(21, 2).Value = "Chainlink"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "7.478"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = "  +0.44%  "
$ws.Cells.Item(22, 2).Value = "BinanceUSD"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "1.001"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = "  +0.06%  "
$ws.Cells.Item(23, 2).Value = "Monero"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "158.06"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = "  +0.39%  "
$ws.Cells.Item(24, 2).Value = "Stellar"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.1382"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = "  -1.13%  "
$ws.Cells.Item(25, 2).Value = "Cosmos"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "8.422"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = "  +0.70%  "
$ws.Cells.Item(26, 2).Value = "EthereumClassic"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "17.71"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = "  +0.46%  "
$ws.Cells.Item(27, 2).Value = "Toncoin"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "1.385"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = "  +6.73%  "
$ws.Cells.Item(28, 2).Value = "PancakeSwap"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.462"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = "  -0.23%  "
$ws.Cells.Item(29, 2).Value = "Hedera"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.05596"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = "  +0.18%  "
$ws.Cells.Item(30, 2).Value = "Filecoin"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "4.133"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = "  +0.74%  "
$ws.Cells.Item(31, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.066"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = "  +0.85%  "
$ws.Cells.Item(32, 2).Value = "LidoDAOToken"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.842"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value = "  -0.11%  "
$ws.Cells.Item(33, 2).Value = "ARBITRUM"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.165"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = "  +0.72%  "
$ws.Cells.Item(34, 2).Value = "ImmutableX"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.6966"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -1.85%  "
$ws.Cells.Item(35, 2).Value = "HuobiToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.588"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = "  +0.09%  "
$ws.Cells.Item(36, 2).Value = "VeChain"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.01804"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = "  +0.19%  "
$ws.Cells.Item(37, 2).Value = "Maker"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(37, 4).Value = "1.229.67"
$ws.Cells.Item(37, 5).Value = "  -0.38%  "
$ws.Cells.Item(38, 2).Value = "MXToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.718"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = "  -1.72%  "
$ws.Cells.Item(39, 2).Value = "FraxShare"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "6.433"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = "  -0.09%  "
$ws.Cells.Item(40, 2).Value = "TrustWalletToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.9073"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = "  +0.25%  "
$ws.Cells.Item(41, 2).Value = "PaxDollar"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.001"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = "  +0.08%  "
$ws.Cells.Item(42, 2).Value = "Quant"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "101.80"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = "  -0.13%  "
$ws.Cells.Item(43, 2).Value = "Aave"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "66.04"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -0.12%  "
$ws.Cells.Item(44, 2).Value = "Aptos"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "7.201"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = "  +0.40%  "
$ws.Cells.Item(45, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.00000000119"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = "  -2.36%  "
$ws.Cells.Item(46, 2).Value = "TheSandbox"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.4017"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = "  -0.13%  "
$ws.Cells.Item(47, 2).Value = "EnergySwap"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "9.029"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = "  +0.74%  "
$ws.Cells.Item(48, 2).Value = "RenderToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.685"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = "  +0.46%  "
$ws.Cells.Item(49, 2).Value = "Algorand"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.1150"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = "  +2.57%  "
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.05702"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = "  -0.13%  "
$ws.Cells.Item(51, 2).Value = "Mantle"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.4632"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = "  +0.15%  "
